# edit.ps1 - applies the diff to 大專生研究計畫.docx
#
# Summary of the change:
#  1) In the Abstract ("(一) 摘要") section, the paragraph that used to read
#     "近年來...問題生成。" is replaced by two new paragraphs:
#       A) two tabs followed by "本研究欲透過「自動問句生成閱讀理解選項」為應用目標，欲透過提"
#       B) (indented, left=480) "升問句生成過程中關鍵字篩選之難度及複雜度，來增進自動問句生成閱讀理解模型對未來教育環境的幫助。"
#          (keeps the _GoBack bookmark at the end, as before)
#     Directly after that paragraph there were two empty paragraphs; now there
#     is only one, so one empty paragraph is removed.
#  2) In the "(二) 研究動機與研究問題" section, the single empty paragraph right
#     after the heading is replaced by two new paragraphs:
#       A) (indented, left=480 firstLine=480) the old abstract sentence, lightly
#          edited ("這個" -> "此", trailing "。" -> "，。")
#       B) a new empty paragraph (carrying rFonts hint=eastAsia formatting)

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Step 1 (do the *later* edit in the document first, so paragraph indices for
# earlier content used in Step 2 stay valid).
#
# "(二) 研究動機與研究問題" heading is followed immediately by a single empty
# paragraph. Locate the heading paragraph by its text, then work on the very
# next paragraph.
# ---------------------------------------------------------------------------
$headingIdx = Find-ParagraphIndex $d "研究動機與研究問題"
if ($headingIdx -eq -1) {
    throw "Could not locate heading paragraph 研究動機與研究問題"
}
$targetPara = $d.Paragraphs.Item($headingIdx + 1)
$targetRange = $targetPara.Range

$motivationXml = "<w:p $wns>" +
    "<w:pPr><w:ind w:left=`"480`" w:firstLine=`"480`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr>" +
    "<w:t>近年來，因為自然語言處理領域的技術較成熟，使得自動問題生成成為開發未來學習系統理想的工具。在此領域中，大多依賴預先選擇「關鍵字」來幫助問題生成，。</w:t>" +
    "</w:r></w:p>" +
    "<w:p $wns><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr></w:p>"

$targetRange.InsertXML($motivationXml)

# ---------------------------------------------------------------------------
# Step 2. The abstract paragraph "近年來...問題生成。" (in the "(一) 摘要"
# section) is split into two new paragraphs.
# ---------------------------------------------------------------------------
$abstractIdx = Find-ParagraphIndex $d "近年來，因為自然語言處理領域的技術較成熟"
if ($abstractIdx -eq -1) {
    throw "Could not locate abstract paragraph"
}
$abstractPara = $d.Paragraphs.Item($abstractIdx)
$abstractParaRange = $abstractPara.Range

$abstractXml = "<w:p $wns>" +
    "<w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr>" +
    "<w:r><w:tab/></w:r>" +
    "<w:r><w:tab/></w:r>" +
    "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr>" +
    "<w:t>本研究欲透過「自動問句生成閱讀理解選項」為應用目標，欲透過提</w:t>" +
    "</w:r></w:p>" +
    "<w:p $wns>" +
    "<w:pPr><w:ind w:left=`"480`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr>" +
    "<w:t>升問句生成過程中關鍵字篩選之難度及複雜度，來增進自動問句生成閱讀理解模型對未來教育環境的幫助。</w:t>" +
    "</w:r>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
    "</w:p>"

$abstractParaRange.InsertXML($abstractXml)

# After the split, the paragraph right after the new pair (paragraph B) is the
# first of the two originally-empty paragraphs; remove it so that two empty
# paragraphs become one (as in the diff).
$abstractIdx2 = Find-ParagraphIndex $d "升問句生成過程中關鍵字篩選之難度及複雜度"
if ($abstractIdx2 -eq -1) {
    throw "Could not locate newly inserted paragraph B"
}
$emptyPara = $d.Paragraphs.Item($abstractIdx2 + 1)
$emptyPara.Range.Delete()

Write-Host "Edit complete."
